$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the Price/Volume columns so Excel
# does not auto-convert numeric-looking strings (e.g. "243.93") into
# real numbers when we assign them through .Value.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "96.786.82"
$ws.Range("E2").Value = "  +1.00%  "

$ws.Range("D3").Value = "3.712.37"
$ws.Range("E3").Value = "  +4.33%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "243.93"
$ws.Range("E5").Value = "  +2.11%  "

$ws.Range("E6").Value = "  +20.19%  "

$ws.Range("D7").Value = "673.61"
$ws.Range("E7").Value = "  +3.19%  "

$ws.Range("E8").Value = "  +6.71%  "

$ws.Range("D9").Value = "1.15"
$ws.Range("E9").Value = "  +9.61%  "

$ws.Range("E10").Value = "  -0.05%  "

$ws.Range("D11").Value = "3.714.70"
$ws.Range("E11").Value = "  +4.47%  "

$ws.Range("D12").Value = "45.66"
$ws.Range("E12").Value = "  +6.24%  "

$ws.Range("E13").Value = "  +2.26%  "

$ws.Range("E14").Value = "  +4.25%  "

$ws.Range("D15").Value = "4.403.26"
$ws.Range("E15").Value = "  +4.34%  "

$ws.Range("D16").Value = "0.0000271"
$ws.Range("E16").Value = "  +5.44%  "

$ws.Range("D17").Value = "96.554.14"
$ws.Range("E17").Value = "  +1.12%  "

$ws.Range("D18").Value = "9.08"
$ws.Range("E18").Value = "  +16.56%  "

$ws.Range("D19").Value = "3.713.78"
$ws.Range("E19").Value = "  +3.96%  "

$ws.Range("D20").Value = "13.07"
$ws.Range("E20").Value = "  +4.15%  "

$ws.Range("D21").Value = "18.70"
$ws.Range("E21").Value = "  +6.17%  "

$ws.Range("D22").Value = "0.548"
$ws.Range("E22").Value = "  +7.04%  "

$ws.Range("D23").Value = "519.69"
$ws.Range("E23").Value = "  +3.64%  "

$ws.Range("E24").Value = "  +2.20%  "

$ws.Range("D25").Value = "0.0000209"
$ws.Range("E25").Value = "  +6.60%  "

$ws.Range("E26").Value = "  +0.60%  "

$ws.Range("D27").Value = "102.27"
$ws.Range("E27").Value = "  +7.42%  "

$ws.Range("D28").Value = "13.19"
$ws.Range("E28").Value = "  +3.98%  "

$ws.Range("D29").Value = "0.170"
$ws.Range("E29").Value = "  +13.40%  "

$ws.Range("D30").Value = "3.10"
$ws.Range("E30").Value = "  +3.87%  "

$ws.Range("D31").Value = "12.21"
$ws.Range("E31").Value = "  +8.12%  "

$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("E33").Value = "  +2.43%  "

$ws.Range("D34").Value = "33.28"
$ws.Range("E34").Value = "  +6.64%  "

$ws.Range("E35").Value = "  +0.83%  "

$ws.Range("D36").Value = "0.599"
$ws.Range("E36").Value = "  +6.82%  "

$ws.Range("D37").Value = "1.74"
$ws.Range("E37").Value = "  +8.18%  "

$ws.Range("E38").Value = "  +0.84%  "

$ws.Range("D39").Value = "614.74"
$ws.Range("E39").Value = "  +0.63%  "

$ws.Range("D40").Value = "42.59"
$ws.Range("E40").Value = "  +25.50%  "

$ws.Range("E41").Value = "  +9.34%  "

$ws.Range("D42").Value = "0.974"
$ws.Range("E42").Value = "  +8.41%  "

$ws.Range("E43").Value = "  +8.73%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").Value = "6.26"
$ws.Range("E45").Value = "  +10.62%  "

$ws.Range("D46").Value = "0.0454"
$ws.Range("E46").Value = "  +8.92%  "

$ws.Range("D47").Value = "0.439"
$ws.Range("E47").Value = "  +28.76%  "

$ws.Range("E48").Value = "  +2.54%  "

$ws.Range("E49").Value = "  +0.35%  "

$ws.Range("D50").Value = "8.69"
$ws.Range("E50").Value = "  +7.28%  "

$ws.Range("D51").Value = "54.69"
$ws.Range("E51").Value = "  +4.68%  "

# Restore the original (default/general) cell formatting so the
# workbook style table is unaffected by the temporary text format.
$ws.Range("D2:E51").ClearFormats()
